$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) values
$ws.Range("D2").Value = '26.786.06'
$ws.Range("D3").Value = '1.646.00'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.76'
$ws.Range("D5").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0628'
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.18'
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0843'
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = '1.642.12'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.528'
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '64.72'
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = '26.784.57'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '214.40'
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.41'
$ws.Range("D20").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '146.18'
$ws.Range("D24").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.15'
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.66'
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0510'
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.37'
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.01'
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Value = '1.289.13'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.536'
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.821'
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.806'
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.23'
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.33'
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = '1.791.08'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '61.74'
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '91.86'
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0₆0104'
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0522'
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.70'
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0971'
$ws.Range("D50").Style = "Normal"

# Update Volume(1h) (column E) values
$ws.Range("E2").Value = '  +0.47%  '
$ws.Range("E3").Value = '  +0.15%  '
$ws.Range("E4").Value = '  +0.21%  '
$ws.Range("E5").Value = '  +0.68%  '
$ws.Range("E6").Value = '  -0.66%  '
$ws.Range("E7").Value = '  +0.31%  '
$ws.Range("E8").Value = '  -0.12%  '
$ws.Range("E9").Value = '  +0.04%  '
$ws.Range("E10").Value = '  -0.56%  '
$ws.Range("E11").Value = '  +0.01%  '
$ws.Range("E12").Value = '  +0.44%  '
$ws.Range("E13").Value = '  -0.82%  '
$ws.Range("E14").Value = '  -0.47%  '
$ws.Range("E15").Value = '  -1.21%  '
$ws.Range("E16").Value = '  +0.24%  '
$ws.Range("E17").Value = '  -1.64%  '
$ws.Range("E18").Value = '  -1.10%  '
$ws.Range("E19").Value = '  +0.32%  '
$ws.Range("E20").Value = '  +0.89%  '
$ws.Range("E21").Value = '  +12.71%  '
$ws.Range("E22").Value = '  -0.94%  '
$ws.Range("E23").Value = '  -2.01%  '
$ws.Range("E24").Value = '  +0.14%  '
$ws.Range("E25").Value = '  +0.27%  '
$ws.Range("E26").Value = '  -1.44%  '
$ws.Range("E27").Value = '  -0.21%  '
$ws.Range("E28").Value = '  -0.73%  '
$ws.Range("E29").Value = '  -1.40%  '
$ws.Range("E30").Value = '  +0.79%  '
$ws.Range("E31").Value = '  -0.68%  '
$ws.Range("E32").Value = '  -1.42%  '
$ws.Range("E33").Value = '  +1.08%  '
$ws.Range("E34").Value = '  -0.25%  '
$ws.Range("E35").Value = '  +1.52%  '
$ws.Range("E36").Value = '  -2.94%  '
$ws.Range("E37").Value = '  +0.49%  '
$ws.Range("E38").Value = '  -1.08%  '
$ws.Range("E39").Value = '  +0.37%  '
$ws.Range("E40").Value = '  -1.19%  '
$ws.Range("E41").Value = '  -0.98%  '
$ws.Range("E42").Value = '  -2.44%  '
$ws.Range("E43").Value = '  +0.47%  '
$ws.Range("E44").Value = '  +3.00%  '
$ws.Range("E45").Value = '  -0.74%  '
$ws.Range("E46").Value = '  +0.56%  '
$ws.Range("E47").Value = '  +0.70%  '
$ws.Range("E48").Value = '  +1.00%  '
$ws.Range("E49").Value = '  -1.18%  '
$ws.Range("E50").Value = '  -0.14%  '
$ws.Range("E51").Value = '  +0.70%  '
